$d = $word.ActiveDocument

# 1. "Μπορεί να συνεχίσει την ροή από το βήμα 2 " -> drop the trailing space
$d.Content.Find.Execute(
    "Μπορεί να συνεχίσει την ροή από το βήμα 2 ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Μπορεί να συνεχίσει την ροή από το βήμα 2", 2)

# 2. "Το σύστημα ειδοποιεί το άτομο/άτομα προσωπικού υπεύθυνα για την έγκριση της άμεσης ειδοποίησης"
#    -> "Το σύστημα ειδοποιεί τους διαχειριστές υπεύθυνους για την έγκριση της άμεσης ειδοποίησης"
#    (unique match on the full original sentence, done before change 5/6 so their shared prefix
#    is still unambiguous for those edits)
$d.Content.Find.Execute(
    "Το σύστημα ειδοποιεί το άτομο/άτομα προσωπικού υπεύθυνα για την έγκριση της άμεσης ειδοποίησης",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Το σύστημα ειδοποιεί τους διαχειριστές υπεύθυνους για την έγκριση της άμεσης ειδοποίησης", 2)

# 3. "βλέπει σε ένα ημερολόγιο ... πότε έχει ραντεβού και με ποιους ..." -> "... πότε έχει επίσκεψη και με ποιους ..."
$d.Content.Find.Execute(
    "βλέπει σε ένα ημερολόγιο τις μέρες και τις ώρες που δουλεύει, που δουλεύει, πότε έχει ραντεβού και με ποιους, και τις μέρες που έχει άδεια",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "βλέπει σε ένα ημερολόγιο τις μέρες και τις ώρες που δουλεύει, που δουλεύει, πότε έχει επίσκεψη και με ποιους, και τις μέρες που έχει άδεια", 2)

# 4. ", φαρμακευτική αγωγή και να δηλώσει ενδιαφέρον περαιτέρω ενημέρωσης " -> drop the trailing space
$d.Content.Find.Execute(
    ", φαρμακευτική αγωγή και να δηλώσει ενδιαφέρον περαιτέρω ενημέρωσης ", $true, $false, $false, $false, $false,
    $true, 1, $false, ", φαρμακευτική αγωγή και να δηλώσει ενδιαφέρον περαιτέρω ενημέρωσης", 2)

# 5. "Το σύστημα ειδοποιεί το άτομο/άτομα προσωπικού υπεύθυνα για την έγκριση της ακύρωσης"
#    -> "Το σύστημα ειδοποιεί τους διαχειριστές υπεύθυνους για την έγκριση της ακύρωσης"
$d.Content.Find.Execute(
    "Το σύστημα ειδοποιεί το άτομο/άτομα προσωπικού υπεύθυνα για την έγκριση της ακύρωσης",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Το σύστημα ειδοποιεί τους διαχειριστές υπεύθυνους για την έγκριση της ακύρωσης", 2)

# 6. "Το σύστημα ειδοποιεί το άτομο/άτομα προσωπικού υπεύθυνα για την έγκριση της " (followed by a
#    separate run "αλλαγής") -> "Το σύστημα ειδοποιεί τους διαχειριστές υπεύθυνους για την έγκριση της "
#    By now this prefix is unique in the document (changes 2 and 5 above already rewrote the other
#    two instances), so this Find only touches the "αλλαγής" paragraph.
$d.Content.Find.Execute(
    "Το σύστημα ειδοποιεί το άτομο/άτομα προσωπικού υπεύθυνα για την έγκριση της ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Το σύστημα ειδοποιεί τους διαχειριστές υπεύθυνους για την έγκριση της ", 2)
